# Add pop up blocker / JS alert / confirmation / prompt test-case rows,
# plus a new "AlertAction" column (K), to the TestData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lay down formatting for the new rows (6,7,8) first, copying the style
#     of row 5 (A:J) which already has the right border/wrap styling. -------
$ws.Range("A5:J5").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$ws.Range("A5:J5").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$ws.Range("A5:J5").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)

# New column K formatting, copied from column J on each existing row.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("J5").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("J5").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("J5").Copy()
$ws.Range("K8").PasteSpecial(-4122)

# --- Now fill in the cell values, in the same order the original author
#     typed them (this drives the shared-string table ordering). -----------
$ws.Range("A6").Value = "TestCase5_validateJSAlerts"
$ws.Range("C6").Value = "kw_alerts"
$ws.Range("A7").Value = "TestCase6_validateJSConfirmation"
$ws.Range("B6").Value = "Validate JavaScript Alerts"
$ws.Range("B7").Value = "Validate JavaScript Confirmations"
$ws.Range("C7").Value = "kw_confirm"
$ws.Range("A8").Value = "TestCase7_validateJSPrompt"
$ws.Range("B8").Value = "Validate JavaScript Prompt"
$ws.Range("C8").Value = "kw_prompt"
$ws.Range("K1").Value = "AlertAction"
$ws.Range("K6").Value = "accept"
$ws.Range("K7").Value = "dismiss"
$ws.Range("K8").Value = "accept"
$ws.Range("I8").Value = "Selenium"
$ws.Range("K5").Value = $ws.Range("J5").Value()

# Row 7 (JS Confirmation) wraps onto two lines like rows 2-4, so it needs the
# taller row height.
$ws.Rows(7).RowHeight = 30

# Column K should be the same width as H:J.
$ws.Columns(11).ColumnWidth = 19

# --- Selection: the author ended up with I1 selected. -----------------------
$ws.Range("I1").Select() | Out-Null
